$wb = $excel.ActiveWorkbook

# Remember the currently active sheet so we can restore it at the end
$originalActiveSheet = $wb.ActiveSheet

# Reference the last existing sheet (2025-12-31) to copy formatting from
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Add new sheet after the last one
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2026-01-07"

# Match sheetPr outline properties used by the other weekly sheets
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# Match page margins used by the other weekly sheets
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "volume"
$ws.Range("D1").Value = "publisher"

# Copy header style (bold, centered, bordered) from an existing sheet's header row
$lastSheet.Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

# Data rows: rank, title, volume, highlighted(1/0 - highlighted rows use the
# same "new series" yellow-fill style already used on other weekly sheets)
$data = @(
    @(1, '片田舎のおっさん、剣聖になる~ただの田舎の剣術師範だったのに、大成した弟子たちが俺を放ってくれない件~', 8, 0),
    @(2, '葬送のフリーレン', 15, 0),
    @(3, '俺だけレベルアップな件', 23, 0),
    @(4, 'ゴブリンスレイヤー', 17, 0),
    @(5, '雑用付与術師が自分の最強に気付くまで(コミック)', 10, 0),
    @(6, 'ダーウィン事変', 10, 0),
    @(7, 'ハニーレモンソーダ', 30, 0),
    @(8, 'ゴブリンスレイヤー外伝:イヤーワン', 14, 0),
    @(9, '太陽よりも眩しい星', 13, 0),
    @(10, 'フットボールネーション', 20, 0),
    @(11, '悪食令嬢と狂血公爵 ~その魔物、私が美味しくいただきます!~', 12, 0),
    @(12, 'ループ7回目の悪役令嬢は、元敵国で自由気ままな花嫁生活を満喫する', 8, 0),
    @(13, 'バトルスタディーズ', 47, 0),
    @(14, '薬屋のひとりごと~猫猫の後宮謎解き手帳~', 21, 0),
    @(15, 'ONE PIECE', 113, 0),
    @(16, 'アイヲンモール異世界店、本日グランドオープン! THE COMIC', 1, 1),
    @(17, 'ガチャを回して仲間を増やす 最強の美少女軍団を作り上げろ THE COMIC', 5, 0),
    @(18, '偽装カレシに愛されてしまいました', 3, 1),
    @(19, '薬屋のひとりごと', 16, 0),
    @(20, 'ジョジョの奇妙な冒険 ザ・ジョジョランズ', 7, 0),
    @(21, 'カラフルグレー', 1, 1),
    @(22, '突然パパになった最強ドラゴンの子育て日記~かわいい娘、ほのぼのと人間界最強に育つ~ THE COMIC', 1, 1),
    @(23, 'ケジメつけさせてもらいます。元ヤン弁護士 東矢斎', 2, 1),
    @(24, '片田舎のおっさん、剣聖になる外伝 はじまりの魔法剣士', 2, 1),
    @(25, '一夜限りのお相手が溺愛先生へと変貌しました1', 1, 1),
    @(26, '国民的アイドルが弟になったら', 5, 0),
    @(27, 'ガチャを回して仲間を増やす 最強の美少女軍団を作り上げろ THE COMIC', 13, 0),
    @(28, '片田舎のおっさん、剣聖になる外伝 竜双剣の軌跡', 2, 1),
    @(29, 'ダイヤモンドの功罪', 9, 0),
    @(30, '隣の元カレくん 単行本版', 5, 0),
    @(31, '八雲さんは餌づけがしたい。 特別読切', 12, 0),
    @(32, '先生のうち、行っていい?ダメでもいくけど。1', 1, 1),
    @(33, 'みいちゃんと山田さん', 5, 0),
    @(34, 'ワンダンス', 15, 0),
    @(35, 'ドンケツ第2章', 16, 0),
    @(36, 'オーイ! とんぼ', 60, 0),
    @(37, '無能と追放された最弱魔法剣士、呪いが解けたので最強へ成り上がる1', 1, 1),
    @(38, '元最強の剣士は、異世界魔法に憧れる THE COMIC', 3, 1),
    @(39, 'ジェネリック彼氏じゃダメですか?~元カレの弟と同居はじめました~', 6, 0),
    @(40, '偽装カレシに愛されてしまいました', 2, 1),
    @(41, '異世界迷宮のオーパーツ', 1, 1),
    @(42, '降り積もれ孤独な死よ', 11, 0),
    @(43, 'クズ旦那と離婚したら、最高の愛を注がれています1', 1, 1),
    @(44, '執事ですがなにか?~幼馴染のパワハラ皇女と絶縁したら、隣国の向日葵王女に拾われたのでこの身を捧げます~1', 1, 1),
    @(45, '修羅幼女の英雄譚~半端者と言われた傭兵、幼女に転生して成り上がる~1', 1, 1),
    @(46, '規格外のダンジョン攻略者、実は異世界帰りの元勇者1', 1, 1),
    @(47, '境界迷宮と異界の魔術師', 1, 1),
    @(48, '人気配信者たちのマネージャーになったら、全員元カノだった 第2話', 2, 1),
    @(49, '乙女ゲームの当て馬悪役令嬢は、王太子殿下の幸せを願います! コミック版', 1, 1),
    @(50, '国民的アイドルが弟になったら', 3, 1),
    @(51, '国民的アイドルが弟になったら', 4, 0),
    @(52, '夜を照らすポラリス~なくした記憶と恋の行方~', 3, 1),
    @(53, '悪魔なボクは退魔師サマに愛されたい!!', 1, 1),
    @(54, '欠けてるふたり~男友達と限界の夜に', 3, 1),
    @(55, '偽装カレシに愛されてしまいました', 6, 0),
    @(56, '異世界迷宮のオーパーツ', 2, 1),
    @(57, '異世界迷宮のオーパーツ', 3, 1),
    @(58, 'ツッコミ待ちの町野さん', 1, 1),
    @(59, 'すみっこ漫画家のトンデモ『裏』事件簿', 1, 1),
    @(60, 'ディグイット', 2, 1),
    @(61, '魔力0で最強の大賢者 ~それは魔法ではない、物理だ!~:', 11, 0),
    @(62, 'ゴブリンスレイヤー外伝2 鍔鳴の太刀《ダイ・カタナ》', 9, 0),
    @(63, '剣聖の幼馴染がパワハラで俺につらく当たるので、絶縁して辺境で魔剣士として出直すことにした。(コミック)', 8, 0),
    @(64, '魔入りました!入間くん', 46, 0),
    @(65, 'お姉ちゃんの翠くん', 10, 0),
    @(66, '私たちはシーツの中で恋をする', 2, 1),
    @(67, '反逆の勇者~スキルを使って腹黒王女のココロとカラダを掌握せよ~', 1, 1),
    @(68, '無能は不要と言われ『時計使い』の僕は職人ギルドから追い出されるも、ダンジョンの深部で真の力に覚醒する THE COMIC', 1, 1),
    @(69, 'ガチャを回して仲間を増やす 最強の美少女軍団を作り上げろ THE COMIC', 4, 0),
    @(70, '元最強の剣士は、異世界魔法に憧れる THE COMIC', 2, 1),
    @(71, 'デレたい彼女の裏表 第2話', 2, 1),
    @(72, '彼の魔族は如何なる魔術をもって防護結界を攻略したか1', 1, 1),
    @(73, '素直になれない雪乙女は眠れる竜騎士に甘くとかされる コミック版', 1, 1),
    @(74, 'こじらせ令嬢の幸せな黒歴史 ~鈍感騎士に溺愛されるための秘密のアプローチ~ コミック版', 1, 1),
    @(75, 'ツッコミ待ちの町野さん', 2, 1),
    @(76, 'ツッコミ待ちの町野さん', 3, 1),
    @(77, '二番目な僕と一番の彼女', 1, 1),
    @(78, 'りゅうとあまがみ', 1, 1),
    @(79, 'きらめきの大和くん☆', 1, 1),
    @(80, 'ねこねこ幼女の愛情ごはん~異世界でもふもふ達に料理を作ります!~', 1, 1),
    @(81, 'あきらめ令嬢は恋心なんていらない。~裏切られたはずなのに、婚約者からの溺愛が止まりません!~', 1, 1),
    @(82, 'ライセット! ~転生令嬢による異世界ハーブアイテム革命~', 1, 1),
    @(83, 'S級ギルドを追放されたけど、実は俺だけドラゴンの言葉がわかるので、気付いたときには竜騎士の頂点を極めてました。', 7, 0),
    @(84, '最凶の魔王に鍛えられた勇者、異世界帰還者たちの学園で無双する', 5, 0),
    @(85, '転生したらスライムだった件', 30, 0),
    @(86, '異世界魔王と召喚少女の奴隷魔術', 28, 0),
    @(87, '攻撃力極振りの最強魔術師~筋力値9999の大剣士、転生して二度目の人生を歩む~', 8, 0),
    @(88, '追放されたS級鑑定士は最強のギルドを創る', 9, 0),
    @(89, '日本へようこそエルフさん。', 12, 0),
    @(90, '勇者パーティーにかわいい子がいたので、告白してみた。(コミック)', 13, 0),
    @(91, '拷問バイトくんの日常', 7, 0),
    @(92, 'スノウボールアース', 10, 0),
    @(93, 'ワールドトリガー', 29, 0),
    @(94, 'アオのハコ', 23, 0),
    @(95, 'SAKAMOTO DAYS', 25, 0),
    @(96, '社内探偵', 67, 0),
    @(97, 'バニシング・ツイン~私の中の君~', 1, 1),
    @(98, '或いは、私の名探偵', 1, 1),
    @(99, '或いは、私の名探偵', 2, 1),
    @(100, '或いは、私の名探偵', 3, 1)
)

$rowIndex = 2
foreach ($entry in $data) {
    $rank = $entry[0]
    $title = $entry[1]
    $volume = $entry[2]
    $styled = $entry[3]

    $ws.Cells.Item($rowIndex, 1).Value = $rank
    $ws.Cells.Item($rowIndex, 2).Value = $title
    $ws.Cells.Item($rowIndex, 3).Value = $volume

    if ($styled -eq 1) {
        $lastSheet.Range("C13").Copy()
        $ws.Cells.Item($rowIndex, 3).PasteSpecial(-4122)
        $ws.Cells.Item($rowIndex, 3).Value = $volume
    }

    $rowIndex++
}

$excel.CutCopyMode = 0
$ws.Range("A1").Select()

# Restore the originally active sheet/tab
$originalActiveSheet.Activate()
